# Generate Report for Handback
#
# Updates the localization-status report after a handback:
#   - flips the status text from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it appears,
#   - stamps the "Latest Handback DateTime" columns with fresh timestamps,
#   - adds the "Latest Target File" / "Latest Handback File" hyperlinked
#     cells (columns F/G) on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$targetMd  = "2c71b1d3-5e1a-4ade-ac15-5c5fd7702db2.md"
$zhXlf     = "2c71b1d3-5e1a-4ade-ac15-5c5fd7702db2.3882b98b65e46b376c7d770165d4870736c6b157.zh-cn.xlf"
$deXlf     = "2c71b1d3-5e1a-4ade-ac15-5c5fd7702db2.3882b98b65e46b376c7d770165d4870736c6b157.de-de.xlf"

$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/13cef08e987d520f17a270827496d25560375bc2/e2e/$targetMd"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ee423ef070a52d92bf2fb24ef7de1b2306100768/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/$zhXlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2a0789aa4497311cf9316a57644040eff9615df9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/$deXlf"

# ---------------------------------------------------------------------------
# 1. Overview sheet: flip the status cells.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: flip status, stamp handback datetime, add F/G hyperlinks.
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

# "Latest Handback DateTime" (column H) refresh.
$zh.Range("H2").Value = "2016-03-14 08:44:56"
$zh.Range("H3").Value = "2016-03-14 08:44:56"

# "Latest Target File" (F) / "Latest Handback File" (G).
$zh.Hyperlinks.Add($zh.Range("F2"), $mdUrl, "", "", $targetMd)
$zh.Hyperlinks.Add($zh.Range("G2"), $zhXlfUrl, "", "", $zhXlf)
$zh.Hyperlinks.Add($zh.Range("F3"), $mdUrl, "", "", $targetMd)
$zh.Hyperlinks.Add($zh.Range("G3"), $zhXlfUrl, "", "", $zhXlf)

# ---------------------------------------------------------------------------
# 3. de-de sheet: flip status, stamp a NEW handback datetime, add F/G links.
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# "Latest Handback DateTime" (column H) refresh - distinct timestamp.
$de.Range("H2").Value = "2016-03-14 08:45:02"
$de.Range("H3").Value = "2016-03-14 08:45:02"

# "Latest Target File" (F) / "Latest Handback File" (G).
$de.Hyperlinks.Add($de.Range("F2"), $mdUrl, "", "", $targetMd)
$de.Hyperlinks.Add($de.Range("G2"), $deXlfUrl, "", "", $deXlf)
$de.Hyperlinks.Add($de.Range("F3"), $mdUrl, "", "", $targetMd)
$de.Hyperlinks.Add($de.Range("G3"), $deXlfUrl, "", "", $deXlf)
